# Add a new "Kaplama Var" option row to the options sheet, splitting the
# former combined "Yok" choice (which carried a "Yok,SA,MA,OA" suffix list)
# into a dedicated Kaplama_Var_Mi_opts entry, and update the now-shifted
# Govde_Sonu_Tipi_opts rows whose PrereqFieldKey/PrereqAllowValues point at
# the ground-connection-type options (they should point at the base field
# key "Yer_Baglanti_tipi" with the narrowed allow-list "SA,MA,OA").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("options")

# Insert a new row right after the existing "Kaplama_Var_Mi_opts / Yok" row
# (row 200), pushing every subsequent row down by one.
$ws.Rows("200:200").Insert()

$ws.Range("A200").Value = "Kaplama_Var_Mi_opts"
$ws.Range("B200").Value = "Var"
$ws.Range("C200").Value = "Kaplama Var"
$ws.Range("D200").Value = 2

# These rows used to be 212-214 before the insert; after the shift they are
# 213-215. Their PrereqFieldKey (E) and PrereqAllowValues (F) need updating.
$ws.Range("E213").Value = "Yer_Baglanti_tipi"
$ws.Range("F213").Value = "SA,MA,OA"

$ws.Range("E214").Value = "Yer_Baglanti_tipi"
$ws.Range("F214").Value = "SA,MA,OA"

$ws.Range("E215").Value = "Yer_Baglanti_tipi"
$ws.Range("F215").Value = "SA,MA,OA"

# Restore the originally selected cell (shifted down by one row as well).
$ws.Range("F205").Select() | Out-Null
